# excel_writer: also include totals for the balance columns
#
# On the "Gesamtergebnis" (totals) sheet, the balance columns
# (Startguthaben/C and Endsaldo/D) of the "Total" row used to be written as
# "N/A" placeholders. They should instead carry an actual numeric total (0
# in this sample, same as every other totals column on that row), so the
# totals row is now fully numeric.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gesamtergebnis")

# Replace the "N/A" placeholders in the Total row with numeric totals for
# the balance columns.
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0

# Bring that sheet to the front and select the cells that were just filled
# in, matching where the user's attention/focus ended up after the edit.
$ws.Activate()
$ws.Range("C3:D3").Select()
